$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 71-72; this shifts former rows 71-77 down to 73-79,
# matching the diff (old rows preserved, two brand-new rows inserted above them).
$ws.Rows("71:72").Insert()

# Populate the two newly inserted rows with their new content.
$ws.Range("A71").Value = "2000-2002"
$ws.Range("B71").Value = "**秋田県衛生科学研究所** <br> [薬剤耐性菌の浸淫実態解明に関する調査研究（平成12年度〜平成14年度）](https://www.pref.akita.lg.jp/uploads/public/archive_0000088274_00/ek2002_47_5_03.pdf) <br>（秋田県衛生科学研究所報, 47, 24~29, 2003）"
$ws.Range("C71").Value = "未登録"

$ws.Range("A72").Value = "2000-2002"
$ws.Range("B72").Value = "**秋田県衛生科学研究所** <br> [_Campylobacter jejuni_ isolated from retail poultry meat, bovine feces and bile, and human diarrheal samples in Japan: Comparison of serotypes and genotypes](https://www.sciencedirect.com/science/article/pii/S0928824405001252/pdfft?md5=742f3c6d98078caf0f80135c4d8bd5ce&pid=1-s2.0-S0928824405001252-main.pdf) <br> (FEMS Immunology and Medical Microbiology, Volume 45, Issue 2, 2005, Pages 311-319)"
$ws.Range("C72").Value = "未登録"
